$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("财联社", "东方财富", "同花顺"),
    @("寒武纪-U", "寒武纪-U", "华胜天成"),
    @("工业富联", "工业富联", "长城军工"),
    @("中国稀土", "岩山科技", "工业富联"),
    @("岩山科技", "华胜天成", "岩山科技"),
    @("宁德时代", "宁德时代", "中芯国际"),
    @("中芯国际", "中国稀土", "寒武纪"),
    @("华胜天成", "长城军工", "吉视传媒"),
    @("长城军工", "东方财富", "中国稀土"),
    @("先导智能", "北方稀土", "比亚迪"),
    @("领益智造", "先导智能", "领益智造"),
    @("国轩高科", "领益智造", "宁德时代"),
    @("吉视传媒", "吉视传媒", "万通发展"),
    @("北方稀土", "三维通信", "北方稀土"),
    @("东方财富", "国轩高科", "通富微电"),
    @("三维通信", "拓维信息", "国轩高科"),
    @("拓维信息", "通富微电", "拓维信息"),
    @("通富微电", "中芯国际", "先导智能"),
    @("英维克", "利欧股份", "三维通信"),
    @("长飞光纤", "歌尔股份", "中国船舶"),
    @("万通发展", "盛和资源", "天融信")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}
